$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 935.5814
$ws.Range("J129").Value = 963.25
$ws.Range("L129").Value = 2889.75
$ws.Range("N129").Value = -12889.75
$ws.Range("H135").Value = 293.76923
$ws.Range("I135").Value = 243.35
$ws.Range("K135").Value = 2190.15
$ws.Range("M135").Value = 344.8499999999999
$ws.Range("H138").Value = 3382.39
$ws.Range("I138").Value = 722.27905
$ws.Range("J138").Value = 5389.14
$ws.Range("K138").Value = 2166.83715
$ws.Range("L138").Value = 16167.42
$ws.Range("M138").Value = 2973.16285
$ws.Range("N138").Value = -26447.42
$ws.Range("H141").Value = 24713.639
$ws.Range("I141").Value = 27229.428
$ws.Range("J141").Value = 3581
$ws.Range("K141").Value = 81688.284
$ws.Range("L141").Value = 10743
$ws.Range("M141").Value = -76508.284
$ws.Range("N141").Value = -21103

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 21211.5
$ws.Range("J23").Value = 21211.5
$ws.Range("L23").Value = 21211.5
$ws.Range("N23").Value = -21729.5
$ws.Range("H32").Value = 5362.1514
$ws.Range("I32").Value = 4136.6294
$ws.Range("J32").Value = 10877
$ws.Range("K32").Value = 4136.6294
$ws.Range("L32").Value = 10877
$ws.Range("M32").Value = -3849.6294
$ws.Range("N32").Value = -11451
$ws.Range("H37").Value = 31946.143
$ws.Range("J37").Value = 32303.834
$ws.Range("L37").Value = 32303.834
$ws.Range("N37").Value = -32849.834
$ws.Range("H44").Value = 34073.5
$ws.Range("J44").Value = 34073.5
$ws.Range("L44").Value = 34073.5
$ws.Range("N44").Value = -35049.5
$ws.Range("H55").Value = 33500.89
$ws.Range("J55").Value = 33500.89
$ws.Range("L55").Value = 33500.89
$ws.Range("N55").Value = -34130.89
$ws.Range("H74").Value = 2998.475
$ws.Range("I74").Value = 3117
$ws.Range("J74").Value = 2590.2222
$ws.Range("K74").Value = 3117
$ws.Range("L74").Value = 2590.2222
$ws.Range("M74").Value = -2243
$ws.Range("N74").Value = -4338.2222
$ws.Range("H77").Value = 2998.475
$ws.Range("I77").Value = 3117
$ws.Range("J77").Value = 2590.2222
$ws.Range("K77").Value = 15585
$ws.Range("L77").Value = 12951.111
$ws.Range("M77").Value = -11217
$ws.Range("N77").Value = -21687.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9241.272000000001
$ws.Range("I20").Value = 3301.0833
$ws.Range("J20").Value = 16369.5
$ws.Range("K20").Value = 3301.0833
$ws.Range("L20").Value = 16369.5
$ws.Range("M20").Value = -3054.0833
$ws.Range("N20").Value = -16863.5
$ws.Range("H86").Value = 1737.7059
$ws.Range("I86").Value = 1484.2727
$ws.Range("J86").Value = 2202.3333
$ws.Range("K86").Value = 1484.2727
$ws.Range("L86").Value = 2202.3333
$ws.Range("M86").Value = -361.2727
$ws.Range("N86").Value = -4448.3333
$ws.Range("H89").Value = 1737.7059
$ws.Range("I89").Value = 1484.2727
$ws.Range("J89").Value = 2202.3333
$ws.Range("K89").Value = 7421.363499999999
$ws.Range("L89").Value = 11011.6665
$ws.Range("M89").Value = -1805.363499999999
$ws.Range("N89").Value = -22243.6665
$ws.Range("H134").Value = 2030.125
$ws.Range("I134").Value = 1369.037
$ws.Range("J134").Value = 5600
$ws.Range("K134").Value = 4107.111
$ws.Range("L134").Value = 16800
$ws.Range("M134").Value = -1572.111
$ws.Range("N134").Value = -21870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9092755
$ws.Range("I31").Value = 1136.1714
$ws.Range("J31").Value = 25003088
$ws.Range("K31").Value = 1136.1714
$ws.Range("L31").Value = 25003088
$ws.Range("M31").Value = -841.1713999999999
$ws.Range("N31").Value = -25003678
$ws.Range("H34").Value = 9092755
$ws.Range("I34").Value = 1136.1714
$ws.Range("J34").Value = 25003088
$ws.Range("K34").Value = 1136.1714
$ws.Range("L34").Value = 25003088
$ws.Range("M34").Value = -934.1713999999999
$ws.Range("N34").Value = -25003492
$ws.Range("H68").Value = 80863.25
$ws.Range("J68").Value = 80863.25
$ws.Range("L68").Value = 80863.25
$ws.Range("N68").Value = -82361.25
$ws.Range("H71").Value = 80863.25
$ws.Range("J71").Value = 80863.25
$ws.Range("L71").Value = 242589.75
$ws.Range("N71").Value = -250077.75
$ws.Range("H132").Value = 4303.5713
$ws.Range("I132").Value = 3736.9473
$ws.Range("J132").Value = 5499.778
$ws.Range("K132").Value = 11210.8419
$ws.Range("L132").Value = 16499.334
$ws.Range("M132").Value = -8680.841899999999
$ws.Range("N132").Value = -21559.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4840
$ws.Range("I56").Value = 4840
$ws.Range("K56").Value = 4840
$ws.Range("M56").Value = -4310
$ws.Range("H113").Value = 627.125
$ws.Range("I113").Value = 617.4706
$ws.Range("J113").Value = 650.5714
$ws.Range("K113").Value = 1852.4118
$ws.Range("L113").Value = 1951.7142
$ws.Range("M113").Value = 317.5882000000001
$ws.Range("N113").Value = -6291.7142
$ws.Range("H137").Value = 2851.4285
$ws.Range("J137").Value = 3822.1428
$ws.Range("L137").Value = 11466.4284
$ws.Range("N137").Value = -21666.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6131
$ws.Range("I70").Value = 5577.5
$ws.Range("J70").Value = 7730
$ws.Range("K70").Value = 5577.5
$ws.Range("L70").Value = 7730
$ws.Range("M70").Value = -5307.5
$ws.Range("N70").Value = -8270
$ws.Range("H73").Value = 6131
$ws.Range("I73").Value = 5577.5
$ws.Range("J73").Value = 7730
$ws.Range("K73").Value = 5577.5
$ws.Range("L73").Value = 7730
$ws.Range("M73").Value = -4641.5
$ws.Range("N73").Value = -9602

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 34523.477
$ws.Range("J47").Value = 34523.477
$ws.Range("L47").Value = 34523.477
$ws.Range("N47").Value = -35503.477
$ws.Range("H52").Value = 34523.477
$ws.Range("J52").Value = 34523.477
$ws.Range("L52").Value = 34523.477
$ws.Range("N52").Value = -34989.477
$ws.Range("H122").Value = 5001.7646
$ws.Range("I122").Value = 3083.75
$ws.Range("J122").Value = 6706.6665
$ws.Range("K122").Value = 9251.25
$ws.Range("L122").Value = 20119.9995
$ws.Range("M122").Value = -6801.25
$ws.Range("N122").Value = -25019.9995
$ws.Range("H132").Value = 3532.1365
$ws.Range("I132").Value = 1503.8276
$ws.Range("J132").Value = 7453.533
$ws.Range("K132").Value = 4511.4828
$ws.Range("L132").Value = 22360.599
$ws.Range("M132").Value = -1981.4828
$ws.Range("N132").Value = -27420.599
$ws.Range("H136").Value = 2523.0908
$ws.Range("I136").Value = 1298.48
$ws.Range("J136").Value = 6350
$ws.Range("K136").Value = 3895.44
$ws.Range("L136").Value = 19050
$ws.Range("M136").Value = -1345.44
$ws.Range("N136").Value = -24150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6537593
$ws.Range("I132").Value = 1016.65717
$ws.Range("J132").Value = 20836354
$ws.Range("K132").Value = 3049.97151
$ws.Range("L132").Value = 62509062
$ws.Range("M132").Value = -519.9715099999999
$ws.Range("N132").Value = -62514122
$ws.Range("H136").Value = 4627.2104
$ws.Range("I136").Value = 1411.5555
$ws.Range("J136").Value = 7521.3
$ws.Range("K136").Value = 4234.666499999999
$ws.Range("L136").Value = 22563.9
$ws.Range("M136").Value = -1684.666499999999
$ws.Range("N136").Value = -27663.9
